$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("BonusPower", 0.01, 0,   7,   100, 0,   "lose"),
    @("BonusPower", 2,    810, 111, 13,  26,  "win"),
    @("BonusPower", 2,    920, 137, 70,  140, "win")
)

$startRow = 36
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($r, $col).Value = $rowData[$c]
    }
}
